# Update the cached regression-table figures on the "fc_robustness" sheet.
# The workbook pulls its numbers from an external linked CSV
# ([1]fc_robustness!...); since that external source can't be refreshed
# here, we write the refreshed figures straight into the display cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = "-204.0***"
    "C4"  = "-299.9***"
    "D4"  = "-207.7***"
    "E4"  = "-98.5***"
    "F4"  = "-146.3**"

    "B5"  = "(48.1)"
    "C5"  = "(83.3)"
    "D5"  = "(49.0)"
    "E5"  = "(36.7)"

    "B6"  = "-38.9"
    "C6"  = "-56.4"
    "D6"  = "-32.6"
    "E6"  = "-30.7"
    "F6"  = "-25.3"

    "E7"  = "(39.2)"
    "F7"  = "(74.4)"

    "B11" = "942.4"
    "C11" = "1389.9"
    "D11" = "1026.1"
    "E11" = "480.7"
    "F11" = "927.7"

    "C15" = "-0.22***"
    "E15" = "-0.062***"
    "F15" = "-0.097**"

    "C16" = "(0.051)"
    "F16" = "(0.044)"

    "B17" = "-0.0086"
    "C17" = "-0.053"
    "D17" = "-0.0035"
    "E17" = "-0.031*"
    "F17" = "-0.043"

    "C18" = "(0.045)"
    "E18" = "(0.018)"

    "D21" = "0.027"

    "C22" = "1.12"
    "D22" = "0.72"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
